$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("#Scopes")
$ws2 = $wb.Worksheets.Item("#TParties")
$ws3 = $wb.Worksheets.Item("#Concerns")
$ws4 = $wb.Worksheets.Item("#TTexts")
$ws5 = $wb.Worksheets.Item("#Organizations")

# -----------------------------------------------------------------
# Sheet "#Concerns" (now holds Objectives content instead of Concerns)
# -----------------------------------------------------------------

# Header row renamed from Concerns-related labels to Objectives-related labels
$ws3.Range("A1").Value = "[Objectives]"
$ws3.Range("B1").Value = "ttIsaObjective"
$ws3.Range("E1").Value = "objvSHRoleName"

# Row 3 - Transportation need
$ws3.Range("A3").Formula = '=IF(OR($C3="",$D3=""),"",CONCATENATE("Obj_",$C3,"_",$D3))'
$ws3.Range("B3").Formula = '=IF($A3="","",$A3)'
$ws3.Range("F3").Value = "A box that [Sender] refers to as '[Parcel]', needs to be transported from [SenderAddress] to [DeliveryAddress]."

# Row 4 - Transportation offer
$ws3.Range("A4").Formula = '=IF(OR($C4="",$D4=""),"",CONCATENATE("Obj_",$C4,"_",$D4))'
$ws3.Range("B4").Formula = '=IF($A4="","",$A4)'
$ws3.Range("F4").Value = "The risks I run when committing myself to ``pick up a parcel (with barcode [Barcode]) at [SenderAddress], and deliver it to [Recipient] at [DeliveryAddress]``, are acceptable."

# Row 5 - Commit to pay
$ws3.Range("A5").Formula = '=IF(OR($C5="",$D5=""),"",CONCATENATE("Obj_",$C5,"_",$D5))'
$ws3.Range("B5").Formula = '=IF($A5="","",$A5)'
$ws3.Range("F5").Value = "[Sender] has committed to pay [TransportationFee] when [Barcode] is delivered to [Recipient]."

# Row 6 - Parcel identifiability by Transporter
$ws3.Range("A6").Formula = '=IF(OR($C6="",$D6=""),"",CONCATENATE("Obj_",$C6,"_",$D6))'
$ws3.Range("B6").Formula = '=IF($A6="","",$A6)'
$ws3.Range("F6").Value = "I have labeled [Parcel] with the barcode ([Barcode]) that [Transporter] has provided me with."

# Row 8 - Parcel delivery
$ws3.Range("A8").Formula = '=IF(OR($C8="",$D8=""),"",CONCATENATE("Obj_",$C8,"_",$D8))'
$ws3.Range("B8").Formula = '=IF($A8="","",$A8)'
$ws3.Range("F8").Value = "[Transporter] needs box with [Barcode] to be delivered at [DeliveryAddress]."

# Row 9 - Recipient signature
$ws3.Range("A9").Formula = '=IF(OR($C9="",$D9=""),"",CONCATENATE("Obj_",$C9,"_",$D9))'
$ws3.Range("B9").Formula = '=IF($A9="","",$A9)'
$ws3.Range("F9").Value = "[ReceiptSignature] is the signature of [Recipient] under [DeliveryReceipt]."

# Column B is no longer hidden
$ws3.Columns.Item(2).Hidden = $false

# -----------------------------------------------------------------
# Sheet "#TTexts"
# -----------------------------------------------------------------

# New header labels in column F
$ws4.Range("F1").Value = "ttReqdLoA"
$ws4.Range("F2").Value = "ISOLevel"

# Descriptions get a trailing period (minor copy-edit pass)
$ws4.Range("E3").Value = "Identifier by which [Sender] identifies the parcel."
$ws4.Range("E4").Value = "Identifier by which [Transporter] can identify the parcel to be transported."
$ws4.Range("E5").Value = "the address/location where the parcel has to be picked up."
$ws4.Range("E6").Value = "the address/location where the parcel has to be delivered."
$ws4.Range("E7").Value = "(Name of) the party that has to receive the parcel."
$ws4.Range("E8").Value = "Fee for transporting a box (size = [Dimensions] cm3, weight = [Weight] grams), from [SenderAddress] to [DeliveryAddress], in Euro's."
$ws4.Range("F8").Value = 2
$ws4.Range("E9").Value = "Dimensions of the parcel (length x width x height, all in cm)."
$ws4.Range("E10").Value = "weight of the parcel (in grams)."
$ws4.Range("E12").Value = "Identifier by which [Transporter] can identify the parcel to be transported."
$ws4.Range("E13").Value = "the address/location where the parcel has to be delivered."
$ws4.Range("E14").Value = "Statement saying that [Recipient] has received a parcel with [Barcode]."
$ws4.Range("E15").Value = "Independently verifiable claim by [Recipient] that ``[DeliveryReceipt]`` is truthful."
$ws4.Range("F15").Value = 2

# -----------------------------------------------------------------
# View / selection updates
# -----------------------------------------------------------------

# #Concerns sheet: selection moves from F8 to F9
$ws3.Activate()
$ws3.Range("F9").Select()

# #TTexts sheet: selection moves from E10 to E15
$ws4.Activate()
$ws4.Range("E15").Select()

# #TParties becomes the active tab when the workbook is opened
$ws2.Activate()
$ws2.Range("B1").Select()
